$d = $word.ActiveDocument

$d.Content.Find.Execute("360÷6=60, 0", $true, $false, $false, $false, $false, $true, 1, $false, "661÷3=220, 1", 2)
$d.Content.Find.Execute("907÷3=302, 1", $true, $false, $false, $false, $false, $true, 1, $false, "143÷8=17, 7", 2)
$d.Content.Find.Execute("151÷7=21, 4", $true, $false, $false, $false, $false, $true, 1, $false, "774÷2=387, 0", 2)
$d.Content.Find.Execute("727÷9=80, 7", $true, $false, $false, $false, $false, $true, 1, $false, "444÷4=111, 0", 2)
$d.Content.Find.Execute("876÷8=109, 4", $true, $false, $false, $false, $false, $true, 1, $false, "310÷8=38, 6", 2)
$d.Content.Find.Execute("445÷3=148, 1", $true, $false, $false, $false, $false, $true, 1, $false, "148÷7=21, 1", 2)
$d.Content.Find.Execute("288÷8=36, 0", $true, $false, $false, $false, $false, $true, 1, $false, "579÷2=289, 1", 2)
$d.Content.Find.Execute("595÷3=198, 1", $true, $false, $false, $false, $false, $true, 1, $false, "799÷5=159, 4", 2)
$d.Content.Find.Execute("603÷6=100, 3", $true, $false, $false, $false, $false, $true, 1, $false, "294÷2=147, 0", 2)
$d.Content.Find.Execute("882÷4=220, 2", $true, $false, $false, $false, $false, $true, 1, $false, "965÷6=160, 5", 2)
$d.Content.Find.Execute("817÷3=272, 1", $true, $false, $false, $false, $false, $true, 1, $false, "993÷6=165, 3", 2)
$d.Content.Find.Execute("147÷4=36, 3", $true, $false, $false, $false, $false, $true, 1, $false, "579÷4=144, 3", 2)
$d.Content.Find.Execute("121÷6=20, 1", $true, $false, $false, $false, $false, $true, 1, $false, "142÷9=15, 7", 2)
$d.Content.Find.Execute("308÷3=102, 2", $true, $false, $false, $false, $false, $true, 1, $false, "707÷4=176, 3", 2)
$d.Content.Find.Execute("325÷5=65, 0", $true, $false, $false, $false, $false, $true, 1, $false, "707÷7=101, 0", 2)
$d.Content.Find.Execute("106÷7=15, 1", $true, $false, $false, $false, $false, $true, 1, $false, "300÷6=50, 0", 2)
$d.Content.Find.Execute("901÷8=112, 5", $true, $false, $false, $false, $false, $true, 1, $false, "121÷7=17, 2", 2)
$d.Content.Find.Execute("222÷9=24, 6", $true, $false, $false, $false, $false, $true, 1, $false, "514÷4=128, 2", 2)
$d.Content.Find.Execute("104÷2=52, 0", $true, $false, $false, $false, $false, $true, 1, $false, "953÷9=105, 8", 2)
$d.Content.Find.Execute("774÷3=258, 0", $true, $false, $false, $false, $false, $true, 1, $false, "690÷3=230, 0", 2)
$d.Content.Find.Execute("938÷7=134, 0", $true, $false, $false, $false, $false, $true, 1, $false, "848÷4=212, 0", 2)
$d.Content.Find.Execute("711÷2=355, 1", $true, $false, $false, $false, $false, $true, 1, $false, "944÷4=236, 0", 2)
$d.Content.Find.Execute("867÷9=96, 3", $true, $false, $false, $false, $false, $true, 1, $false, "189÷5=37, 4", 2)
$d.Content.Find.Execute("975÷3=325, 0", $true, $false, $false, $false, $false, $true, 1, $false, "623÷2=311, 1", 2)
$d.Content.Find.Execute("534÷8=66, 6", $true, $false, $false, $false, $false, $true, 1, $false, "831÷3=277, 0", 2)
